# Update the "取得日時" (acquired datetime) timestamp in rows 2-8 of the
# "ランサーズ" sheet from 2025-09-06 18:20:50 to 2025-09-06 18:28:51.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-09-06 18:28:51"
}
